$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3828547996091605
$ws.Range("D2").Value = 0.7055015595029293

$ws.Range("C3").Value = -0.5629727912688752
$ws.Range("D3").Value = 0.5791462298687748

$ws.Range("C4").Value = 0.6184602870273573
$ws.Range("D4").Value = 0.5426228001449704

$ws.Range("C5").Value = -0.1047573264290098
$ws.Range("D5").Value = 0.9175175942500626

$ws.Range("C6").Value = -0.7305496900641796
$ws.Range("D6").Value = 0.4727631458341104

$ws.Range("C7").Value = 0.1214704558035892
$ws.Range("D7").Value = 0.9044210904474679

$ws.Range("C8").Value = -0.3781875218724352
$ws.Range("D8").Value = 0.7089158956829102

$ws.Range("C9").Value = 0.8955111711588992
$ws.Range("D9").Value = 0.3802078991838989

$ws.Range("C10").Value = 0.4400260438982854
$ws.Range("D10").Value = 0.6642109374961542

$ws.Range("C11").Value = -0.5862918235365201
$ws.Range("D11").Value = 0.563647147114211
